# Refresh the charging-station report with the latest polling snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data row (2-53) now carries the new "last refreshed" timestamp in column D.
$newTimestamp = 45987.298032407409
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

# Rows 19-53 reflect the refreshed list of stations/terminals/timestamps (A/B/C).
$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "603号直流"
$ws.Cells.Item(19, 3).Value = 45980.250173611108
$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20, 2).Value = "502号直流"
$ws.Cells.Item(20, 3).Value = 45982.555462962962
$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "904号直流"
$ws.Cells.Item(21, 3).Value = 45985.569664351853
$ws.Cells.Item(22, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22, 2).Value = "604号直流"
$ws.Cells.Item(22, 3).Value = 45985.570324074077
$ws.Cells.Item(23, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(23, 2).Value = "201号直流"
$ws.Cells.Item(23, 3).Value = 45985.859444444446
$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(24, 2).Value = "202号直流"
$ws.Cells.Item(24, 3).Value = 45986.069837962961
$ws.Cells.Item(25, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(25, 2).Value = "002A号直流"
$ws.Cells.Item(25, 3).Value = 45986.180451388886
$ws.Cells.Item(26, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(26, 2).Value = "501号直流"
$ws.Cells.Item(26, 3).Value = 45986.210601851853
$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(27, 2).Value = "401号直流"
$ws.Cells.Item(27, 3).Value = 45986.211840277778
$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(28, 2).Value = "011A号直流"
$ws.Cells.Item(28, 3).Value = 45986.25136574074
$ws.Cells.Item(29, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(29, 2).Value = "103号直流"
$ws.Cells.Item(29, 3).Value = 45986.260798611111
$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(30, 2).Value = "103号直流"
$ws.Cells.Item(30, 3).Value = 45986.329641203702
$ws.Cells.Item(31, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value = "805号直流"
$ws.Cells.Item(31, 3).Value = 45986.388541666667
$ws.Cells.Item(32, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(32, 2).Value = "401号直流"
$ws.Cells.Item(32, 3).Value = 45986.421736111108
$ws.Cells.Item(33, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(33, 2).Value = "802号直流"
$ws.Cells.Item(33, 3).Value = 45986.517199074071
$ws.Cells.Item(34, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(34, 2).Value = "305号直流"
$ws.Cells.Item(34, 3).Value = 45986.536064814813
$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(35, 2).Value = "102号直流"
$ws.Cells.Item(35, 3).Value = 45986.537812499999
$ws.Cells.Item(36, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(36, 2).Value = "903号直流"
$ws.Cells.Item(36, 3).Value = 45986.547430555554
$ws.Cells.Item(37, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(37, 2).Value = "801号直流"
$ws.Cells.Item(37, 3).Value = 45986.548611111109
$ws.Cells.Item(38, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(38, 2).Value = "505号直流"
$ws.Cells.Item(38, 3).Value = 45986.555810185186
$ws.Cells.Item(39, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(39, 2).Value = "204号直流"
$ws.Cells.Item(39, 3).Value = 45986.556134259263
$ws.Cells.Item(40, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(40, 2).Value = "107号直流"
$ws.Cells.Item(40, 3).Value = 45986.558680555558
$ws.Cells.Item(41, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(41, 2).Value = "702号直流"
$ws.Cells.Item(41, 3).Value = 45986.565381944441
$ws.Cells.Item(42, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(42, 2).Value = "402号直流"
$ws.Cells.Item(42, 3).Value = 45986.573229166665
$ws.Cells.Item(43, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(43, 2).Value = "503号直流"
$ws.Cells.Item(43, 3).Value = 45986.577627314815
$ws.Cells.Item(44, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(44, 2).Value = "203号直流"
$ws.Cells.Item(44, 3).Value = 45986.586423611108
$ws.Cells.Item(45, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(45, 2).Value = "406号直流"
$ws.Cells.Item(45, 3).Value = 45986.586550925924
$ws.Cells.Item(46, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(46, 2).Value = "101号直流"
$ws.Cells.Item(46, 3).Value = 45986.593726851854
$ws.Cells.Item(47, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(47, 2).Value = "B02号直流"
$ws.Cells.Item(47, 3).Value = 45986.650520833333
$ws.Cells.Item(48, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(48, 2).Value = "A01号直流"
$ws.Cells.Item(48, 3).Value = 45986.706724537034
$ws.Cells.Item(49, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(49, 2).Value = "311号直流"
$ws.Cells.Item(49, 3).Value = 45986.716481481482
$ws.Cells.Item(50, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(50, 2).Value = "301号直流"
$ws.Cells.Item(50, 3).Value = 45986.717476851853
$ws.Cells.Item(51, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(51, 2).Value = "A03号直流"
$ws.Cells.Item(51, 3).Value = 45986.730462962965
$ws.Cells.Item(52, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(52, 2).Value = "804号直流"
$ws.Cells.Item(52, 3).Value = 45986.76394675926
$ws.Cells.Item(53, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(53, 2).Value = "404号直流"
$ws.Cells.Item(53, 3).Value = 45986.769409722219

# Rows 54-57 no longer have entries in the refreshed report.
$ws.Range("A54:D57").ClearContents()

# Restore the author's last selection.
$ws.Range("F12").Select()
